$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the sum formula in D11
$ws.Range("D11").Formula = "=SUM(D2:D8)"

# Apply a yellow fill to the new cell to match style
$ws.Range("D11").Interior.Color = 65535

# Update selection to reflect the new active cell
$ws.Range("C11").Select()
